# Add an "{abstract}" line above "Command" in the SaveCommandUML.png
# diagram on slide 1, and grow the rectangle upward/taller to fit it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Rectangle 8" shape whose text is "Command" (the 4th shape
# on the slide, shape Id 10).
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -eq "Command") {
            $sh = $candidate
            break
        }
    }
}

# Resize/reposition the rectangle: keep left/width, move top up and
# increase the height so the new line fits.
$sh.Top = 90.0
$sh.Height = 41.0471

# Insert a new first paragraph reading "{abstract}" before "Command",
# matching the existing run's bold/italic/size/colour formatting.
$tr = $sh.TextFrame.TextRange
$newRange = $tr.InsertBefore("{abstract}`r")
$newRange.Font.Size = 14
$newRange.Font.Bold = $true
$newRange.Font.Italic = $true
$newRange.ParagraphFormat.Alignment = 2
